# Remove stray trailing spaces from a handful of lyric lines (unicode /
# text-formatting cleanup) without disturbing anything else in the shape
# (run/line-break structure, autofit size, position, etc.).

$p = $ppt.ActivePresentation

function Trim-TrailingSpace($slideIndex, $shapeIndex, $position) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item($shapeIndex)

    # Remember the shape's geometry so we can restore it: editing the
    # TextRange re-triggers the shape's "resize to fit text" (spAutoFit)
    # behaviour, which would otherwise shift the stored height/position.
    $origLeft   = $shape.Left
    $origTop    = $shape.Top
    $origWidth  = $shape.Width
    $origHeight = $shape.Height

    $textRange = $shape.TextFrame.TextRange
    $char = $textRange.Characters($position, 1)

    if ($char.Text -eq " ") {
        $char.Text = ""
    }

    $shape.Left   = $origLeft
    $shape.Top    = $origTop
    $shape.Width  = $origWidth
    $shape.Height = $origHeight
}

# Slide 1, TextBox 1: "Blessed assurance, Jesus is mine! " / "O what a
# foretaste of glory divine! "
# (process highest character index first so that removing one space does
# not shift the index of the other edit within the same shape)
Trim-TrailingSpace 1 1 72
Trim-TrailingSpace 1 1 35

# Slide 3, TextBox 1: "Angels descending bring from above "
Trim-TrailingSpace 3 1 113

# Slide 4, TextBox 1: "Perfect submission, all is at rest "
Trim-TrailingSpace 4 1 36
